$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C numeric values (rows 2-18, skipping unchanged rows)
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 6
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 6
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 7
$ws.Range("C13").Value = 10
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 8
$ws.Range("C18").Value = 6

# Update B13 text label from <delta> to <alpha>
$ws.Range("B13").Value = "<alpha>"
